$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 38
$rng = $ws.Range("A$row" + ":D$row")

# Force literal text entry (prevents Excel from auto-converting the
# date-looking / numeric-looking strings into a date serial or a number),
# then strip the resulting "quote prefix" cell style back to the sheet's
# plain default style so the cells end up as plain shared-string cells,
# matching every other data row already in the sheet.
$ws.Cells.Item($row, 1).Value = "'2026-02-07"
$ws.Cells.Item($row, 2).Value = "'151050"
$ws.Cells.Item($row, 3).Value = "'5"
$ws.Cells.Item($row, 4).Value = "'1"
$rng.Style = "Normal"
